# Add files via upload
# Appends 8 new Q/A rows (rows 593-600) to the single-sheet "input"/"response"
# dataset, reusing existing shared-string text that is already present
# elsewhere in the workbook, then updates the view/selection state to match
# where Excel would land after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each pair below is the exact text already used elsewhere in the sheet, so
# Excel's shared-strings table will reuse the existing <si> entries instead
# of minting new ones.
$newRows = @(
    @('토끼는 왜 시계를 보고 있는 거야?', '토끼는 항상 바빠서 시계를 보고 다니는 것 같아요. 마치 중요한 일이 있나 봐요!'),
    @('여기서 만난 가장 특이한 동물은 뭐야?', '토끼와 체셔 고양이도 특이하지만, 말하는 애벌레도 굉장히 이상하고 신기했어요!'),
    @('우리는 왜 때로 피곤함을 느낄까요?', '피곤함은 몸과 마음이 휴식이 필요하다는 신호예요. 그 신호를 듣고 잠시 쉬는 것이 중요하죠.'),
    @('왜 우리는 자주 피곤함을 느낄까요?', '피로는 과도한 업무나 스트레스, 수면 부족 등으로 발생할 수 있으며, 휴식과 적절한 수면을 취하는 것이 중요해요.'),
    @('어떤 친구가 좋은 친구인가요?', '함께 있을 때 가장 나다운 모습을 보여줄 수 있는 친구가 좋은 친구예요!'),
    @('우리는 왜 행복을 기다리게 될까요?', '행복은 기다리는 것만으로 오는 것이 아니라, 자신이 적극적으로 찾고 만들어가는 것이 중요해요.'),
    @('나는 언제쯤 행복해질까요?', '행복은 기다리는 것이 아니라, 지금 이 순간 속에서 찾는 거예요! 작은 것부터 하나씩 발견해 보세요.'),
    @('어떻게 더 행복해질 수 있을까요?', '행복은 내면의 평화와 삶의 의미를 찾는 과정에서 옵니다. 작은 기쁨을 찾고, 감사하는 마음을 가지면 더 행복할 수 있어요.')
)

$startRow = 593
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $newRows[$i][0]
    $ws.Range("B$row").Value = $newRows[$i][1]
}

$lastRow = $startRow + $newRows.Count - 1

# Match the post-edit view state: the window has scrolled down and the new
# active cell/selection sits just past the last populated row.
[void]$ws.Range("A" + ($lastRow + 1)).Select()
$excel.ActiveWindow.ScrollRow = 583
